$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-10 Tuesday" "2024-12-11 Wednesday"

Replace-Text "699÷7=" "197÷5="
Replace-Text "582÷3=" "649÷5="
Replace-Text "527÷2=" "635÷2="
Replace-Text "996÷5=" "727÷3="
Replace-Text "121÷8=" "551÷9="

Replace-Text "707÷3=" "854÷4="
Replace-Text "611÷2=" "850÷7="
Replace-Text "990÷9=" "557÷7="
Replace-Text "488÷9=" "246÷6="
Replace-Text "986÷8=" "586÷3="

Replace-Text "308÷5=" "255÷6="
Replace-Text "312÷3=" "909÷3="
Replace-Text "838÷9=" "524÷4="
Replace-Text "902÷5=" "604÷5="
Replace-Text "145÷8=" "808÷9="

Replace-Text "905÷5=" "790÷9="
Replace-Text "628÷2=" "886÷8="
Replace-Text "880÷7=" "628÷7="
Replace-Text "589÷2=" "436÷2="
Replace-Text "947÷2=" "449÷5="

Replace-Text "428÷6=" "914÷5="
Replace-Text "508÷6=" "991÷3="
Replace-Text "529÷6=" "449÷4="
Replace-Text "398÷5=" "724÷6="
Replace-Text "295÷8=" "382÷8="
